$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Taylor Wessing", "18sec", "1"),
    @("Samvad Partners", "3sec", "0"),
    @("Remfry And Sagar", "22sec", "0"),
    @("Conyers", "26sec", "1"),
    @("Krogerus", "8sec", "1"),
    @("Zepos And Yannopoulos", "25sec", "0"),
    @("Fangda Partners", "9sec", "0"),
    @("JSA", "16sec", "1"),
    @("TC Law Firm", "14sec", "0"),
    @("Clemens Law", "6sec", "1"),
    @("Paul Hastings", "11sec", "1"),
    @("Carey Olsen", "27sec", "1"),
    @("Dentons", "19sec", "0"),
    @("Legance", "25sec", "0"),
    @("EBN", "6sec", "0"),
    @("Appleby Global", "25sec", "1"),
    @("Higgs And Johnson", "7sec", "1"),
    @("Walkers", "43sec", "1"),
    @("Ecija", "1min 3sec", "0"),
    @("DahlLaw", "2min 35sec", "1"),
    @("Gide Loyrette Nouel", "14sec", "1"),
    @("Gornitzky And Co", "10sec", "0"),
    @("GÖRG", "39sec", "0"),
    @("Beauchamps", "6sec", "0"),
    @("Lex Caribbean", "4sec", "0"),
    @("Portolano Cavallo", "8sec", "0"),
    @("Ashurst", "22sec", "1"),
    @("Cobalt Legal", "8sec", "0"),
    @("Pavia And Ansaldo", "9sec", "0"),
    @("KRB Law Firm", "14sec", "1"),
    @("Dechert LLP", "6sec", "1"),
    @("DSK Legal", "12sec", "1"),
    @("Spruson And Ferguson", "9sec", "1"),
    @("Dillon Eustace", "8sec", "1"),
    @("White and Case", "23sec", "1"),
    @("Veritas Legal", "13sec", "1"),
    @("BonelliErede", "27sec", "1"),
    @("Harneys", "2min 2sec", "1"),
    @("Meitar Law Offices", "21sec", "1"),
    @("Al Tamimi", "23sec", "1"),
    @("Consortium Legal", "15sec", "1"),
    @("Anand And Anand", "13sec", "1"),
    @("BNT", "15sec", "1"),
    @("Hankun Law", "13sec", "1"),
    @("Mourant", "7sec", "0"),
    @("AL Goodbody", "9sec", "0"),
    @("Carnelutti Law Firm", "6sec", "0"),
    @("Tanner DeWitt", "9sec", "0"),
    @("Pedersoli", "11sec", "1"),
    @("Grandall", "47sec", "0"),
    @("Kromann Reumert", "16sec", "1")
)

# Column C holds digit-only strings ("0"/"1") that must remain text rather
# than be auto-coerced to numbers by Excel's type inference. Temporarily
# mark the range as Text, write the values, then restore the General
# format so the cell style matches the rest of the sheet.
$colC = $ws.Range("C2:C52")
$colC.NumberFormat = "@"

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$colC.NumberFormat = "General"
